$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mods")

# New header for column C
$ws.Range("C1").Value = "power"
$ws.Range("C1").Style = "Nadpis 2"

# New power values per row (matching vendor/name rows 2-6)
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 25
$ws.Range("C4").Value = 70
$ws.Range("C5").Value = 50
$ws.Range("C6").Value = 80

# Set column width for C to match the diff (closest attainable snap value)
$ws.Columns.Item(3).ColumnWidth = 18.5

# Update selection to match the diff (activeCell C6)
$ws.Range("C6").Select()
